$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '28.188.69'
Set-TextValue "E2" '  -0.96%  '

Set-TextValue "D3" '1.802.57'
Set-TextValue "E3" '  -1.24%  '

Set-TextValue "E4" '  +0.33%  '

Set-TextValue "D5" '315.25'
Set-TextValue "E5" '  +0.22%  '

Set-TextValue "E6" '  +0.29%  '

Set-TextValue "D7" '0.5250'
Set-TextValue "E7" '  +3.09%  '

Set-TextValue "D8" '0.3807'

Set-TextValue "D9" '0.08016'
Set-TextValue "E9" '  +4.39%  '

Set-TextValue "D10" '41.39'
Set-TextValue "E10" '  -1.20%  '

Set-TextValue "D11" '1.096'
Set-TextValue "E11" '  -1.30%  '

Set-TextValue "D12" '6.332'
Set-TextValue "E12" '  +0.97%  '

Set-TextValue "D13" '1.004'
Set-TextValue "E13" '  +0.32%  '

Set-TextValue "D14" '20.60'
Set-TextValue "E14" '  -1.98%  '

Set-TextValue "B15" 'WrappedEther'
Set-TextValue "C15" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue "D15" '1.811.96'
Set-TextValue "E15" '  -0.79%  '

Set-TextValue "B16" 'Chainlink'
Set-TextValue "C16" 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D16" '7.335'
Set-TextValue "E16" '  -2.56%  '

Set-TextValue "D17" '92.96'
Set-TextValue "E17" '  -0.05%  '

Set-TextValue "D18" '0.00001092'
Set-TextValue "E18" '  -1.22%  '

Set-TextValue "D19" '0.06613'
Set-TextValue "E19" '  -0.41%  '

Set-TextValue "E20" '  +0.37%  '

Set-TextValue "D21" '17.37'
Set-TextValue "E21" '  -2.30%  '

Set-TextValue "D22" '5.968'
Set-TextValue "E22" '  -2.45%  '

Set-TextValue "D23" '28.246.15'
Set-TextValue "E23" '  -0.80%  '

Set-TextValue "D24" '11.17'
Set-TextValue "E24" '  -0.50%  '

Set-TextValue "D25" '2.239'
Set-TextValue "E25" '  -0.64%  '

Set-TextValue "D26" '157.59'
Set-TextValue "E26" '  +1.03%  '

Set-TextValue "D27" '20.53'
Set-TextValue "E27" '  -4.64%  '

Set-TextValue "D28" '2.007.99'
Set-TextValue "E28" '  -1.28%  '

Set-TextValue "D29" '2.392'
Set-TextValue "E29" '  -0.47%  '

Set-TextValue "D30" '123.00'
Set-TextValue "E30" '  -1.28%  '

Set-TextValue "D31" '0.1097'
Set-TextValue "E31" '  -0.21%  '

Set-TextValue "D32" '1.057'
Set-TextValue "E32" '  -4.83%  '

Set-TextValue "D33" '3.667'
Set-TextValue "E33" '  +0.20%  '

Set-TextValue "D34" '5.557'
Set-TextValue "E34" '  -1.89%  '

Set-TextValue "D35" '0.07260'
Set-TextValue "E35" '  +3.01%  '

Set-TextValue "D36" '12.15'
Set-TextValue "E36" '  +8.45%  '

Set-TextValue "D37" '0.2161'
Set-TextValue "E37" '  -2.36%  '

Set-TextValue "B38" 'FraxShare'
Set-TextValue "C38" 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D38" '8.822'
Set-TextValue "E38" '  +0.59%  '

Set-TextValue "B39" 'VeChain'
Set-TextValue "C39" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D39" '0.02306'
Set-TextValue "E39" '  -0.90%  '

Set-TextValue "D40" '5.046'
Set-TextValue "E40" '  -2.44%  '

Set-TextValue "D41" '0.6185'
Set-TextValue "E41" '  -1.36%  '

Set-TextValue "D42" '1.163'
Set-TextValue "E42" '  -1.07%  '

Set-TextValue "D43" '1.372'
Set-TextValue "E43" '  -1.28%  '

Set-TextValue "D44" '0.6021'
Set-TextValue "E44" '  +2.32%  '

Set-TextValue "D45" '13.16'
Set-TextValue "E45" '  -2.12%  '

Set-TextValue "D46" '3.774'
Set-TextValue "E46" '  +1.19%  '

Set-TextValue "D47" '126.28'
Set-TextValue "E47" '  +1.61%  '

Set-TextValue "D48" '1.200'
Set-TextValue "E48" '  +0.61%  '

Set-TextValue "D49" '1.929'
Set-TextValue "E49" '  -2.71%  '

Set-TextValue "D50" '0.06829'
Set-TextValue "E50" '  -1.03%  '

Set-TextValue "D51" '72.75'
Set-TextValue "E51" '  -2.10%  '
